# Rotate the "Recorded By" (column G) list left by one entry for every row
# that has multiple comma-separated recorders (i.e. moves the first name to
# the end of the list), across every worksheet in the workbook.

$wb = $excel.ActiveWorkbook

foreach ($ws in $wb.Worksheets) {
    $usedRange = $ws.UsedRange
    $rowCount = $usedRange.Rows.Count
    $startRow = $usedRange.Row

    for ($i = 0; $i -lt $rowCount; $i++) {
        $r = $startRow + $i
        $cell = $ws.Cells.Item($r, 7)  # Column G
        $val = $cell.Value2

        if ($val -ne $null -and $val.ToString().Contains(",")) {
            $parts = $val.ToString().Split(",")
            for ($p = 0; $p -lt $parts.Length; $p++) {
                $parts[$p] = $parts[$p].Trim()
            }
            if ($parts.Length -gt 1) {
                $rotated = ($parts[1..($parts.Length - 1)] + $parts[0]) -join ", "
                $cell.Value2 = $rotated
            }
        }
    }
}
